$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on price cells whose new values would otherwise be auto-converted to numbers
$forceTextRows = @(5,6,7,8,10,11,12,16,18,19,20,21,22,24,25,27,28,29,31,32,33,34,35,36,37,38,39,40,41,43,44,45,46,47,48,49,50,51)
foreach ($r in $forceTextRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range('D2').Value = '60.560.88'
$ws.Range('E2').Value = '  -4.34%  '

$ws.Range('D3').Value = '2.917.46'
$ws.Range('E3').Value = '  -3.64%  '

$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('D5').Value = '530.53'
$ws.Range('E5').Value = '  -5.02%  '

$ws.Range('D6').Value = '145.67'
$ws.Range('E6').Value = '  -6.32%  '

$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.24%  '

$ws.Range('D8').Value = '0.560'
$ws.Range('E8').Value = '  -0.05%  '

$ws.Range('D9').Value = '2.914.94'
$ws.Range('E9').Value = '  -4.08%  '

$ws.Range('D10').Value = '0.110'
$ws.Range('E10').Value = '  -2.99%  '

$ws.Range('D11').Value = '5.94'
$ws.Range('E11').Value = '  -7.34%  '

$ws.Range('D12').Value = '0.357'
$ws.Range('E12').Value = '  -2.82%  '

$ws.Range('D13').Value = '3.436.36'
$ws.Range('E13').Value = '  -3.38%  '

$ws.Range('E14').Value = '  +1.50%  '

$ws.Range('D15').Value = '60.747.81'
$ws.Range('E15').Value = '  -4.08%  '

$ws.Range('D16').Value = '23.07'
$ws.Range('E16').Value = '  -4.56%  '

$ws.Range('D17').Value = '2.944.98'
$ws.Range('E17').Value = '  -2.67%  '

$ws.Range('D18').Value = '0.0000142'
$ws.Range('E18').Value = '  -5.85%  '

$ws.Range('D19').Value = '5.02'
$ws.Range('E19').Value = '  -1.91%  '

$ws.Range('D20').Value = '11.71'
$ws.Range('E20').Value = '  -2.89%  '

$ws.Range('D21').Value = '365.67'
$ws.Range('E21').Value = '  -8.50%  '

$ws.Range('D22').Value = '6.52'
$ws.Range('E22').Value = '  -2.47%  '

$ws.Range('E23').Value = '  -0.10%  '

$ws.Range('D24').Value = '5.66'
$ws.Range('E24').Value = '  -2.20%  '

$ws.Range('D25').Value = '64.38'
$ws.Range('E25').Value = '  -1.76%  '

$ws.Range('D26').Value = '3.056.60'
$ws.Range('E26').Value = '  -3.27%  '

$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').Value = '0.185'
$ws.Range('E27').Value = '  -2.70%  '

$ws.Range('B28').Value = 'Polygon'
$ws.Range('C28').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D28').Value = '0.455'
$ws.Range('E28').Value = '  -2.14%  '

$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.44%  '

$ws.Range('D30').Value = '0.0₃0879'
$ws.Range('E30').Value = '  -10.87%  '

$ws.Range('D31').Value = '7.74'
$ws.Range('E31').Value = '  -11.51%  '

$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  +0.07%  '

$ws.Range('D33').Value = '1.68'
$ws.Range('E33').Value = '  -4.42%  '

$ws.Range('D34').Value = '19.86'
$ws.Range('E34').Value = '  -2.87%  '

$ws.Range('D35').Value = '158.39'
$ws.Range('E35').Value = '  -2.78%  '

$ws.Range('D36').Value = '4.44'
$ws.Range('E36').Value = '  -6.27%  '

$ws.Range('D37').Value = '5.67'
$ws.Range('E37').Value = '  -6.39%  '

$ws.Range('D38').Value = '1.01'
$ws.Range('E38').Value = '  -9.35%  '

$ws.Range('D39').Value = '1.22'
$ws.Range('E39').Value = '  -6.74%  '

$ws.Range('D40').Value = '38.01'
$ws.Range('E40').Value = '  +0.40%  '

$ws.Range('D41').Value = '1.50'
$ws.Range('E41').Value = '  -6.20%  '

$ws.Range('D42').Value = '2.364.41'
$ws.Range('E42').Value = '  -7.12%  '

$ws.Range('D43').Value = '3.75'
$ws.Range('E43').Value = '  -5.52%  '

$ws.Range('D44').Value = '0.647'
$ws.Range('E44').Value = '  -3.43%  '

$ws.Range('D45').Value = '21.13'
$ws.Range('E45').Value = '  -7.81%  '

$ws.Range('D46').Value = '0.0576'
$ws.Range('E46').Value = '  -4.17%  '

$ws.Range('B47').Value = 'FirstDigitalUSD'
$ws.Range('C47').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D47').Value = '0.999'
$ws.Range('E47').Value = '  +0.11%  '

$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').Value = '4.99'
$ws.Range('E48').Value = '  -2.68%  '

$ws.Range('D49').Value = '0.0235'
$ws.Range('E49').Value = '  -6.23%  '

$ws.Range('B50').Value = 'WhiteBITCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D50').Value = '10.37'
$ws.Range('E50').Value = '  -1.06%  '

$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').Value = '0.0933'
$ws.Range('E51').Value = '  -1.15%  '
